# Automatically run the main business logic: append the latest week's
# raw finishing-order data to the rawdata_Clio sheet, then leave that
# sheet active/selected on the new row, ready to print.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rawdata_Clio")

# Make this the active sheet first (matches the saved workbook state).
$ws.Activate()

# New row of raw data (row 4): Week number in A, finishing order for each
# driver in C..H. Column B (date) is intentionally left blank — it only
# gets a value once the race for that week has actually been run.
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 3).Value = 2
$ws.Cells.Item(4, 4).Value = 5
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 3
$ws.Cells.Item(4, 7).Value = 4
$ws.Cells.Item(4, 8).Value = 6

# Select the freshly-added row, matching the workbook's saved selection.
$ws.Range("A4:H4").Select()
